$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 16; $row++) {
    $ws.Cells.Item($row, 3).Value = Get-Date -Year 2023 -Month 9 -Day 6 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
}
